# auto open job when level up
# Insert a new "job change unlocked" level row (Id=3, Level=5) into the
# LevelInfo table, shifting the existing rows down by one, and append a
# new "new job available" row (Id=101) at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Final data for rows 4..13 (columns A:E) after the edit.
# Row 6 is the newly inserted "job change" unlock row; everything that
# used to live in rows 6..11 now lives one row further down (7..12); row
# 13 is the brand-new "new job available" row appended at the end.
# ---------------------------------------------------------------------

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = '|您可以|Red|编辑你的卡片||了，选择你喜欢的卡片痛痛快快的战斗吧！'
$ws.Range("E4").Value = 1

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = '|您可以和您结交的朋友进行|Blue|切磋||了！战胜他们可以获得|Gold|经验，金钱和卡片||哦！'
$ws.Range("E5").Value = 2

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = '|您可以通过城堡面板，转职按钮来完成|Red|转职||，可以随时切换自己的职业。'
$ws.Range("E6").Value = 11

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = '|您可以通过配方合成|Blue|装备||了！装备可以使你在战斗中更有优势！'
$ws.Range("E7").Value = 3

$ws.Range("A8").Value = 5
$ws.Range("B8").Value = 10
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = '|您可以使用|Gold|商城||了。通过商城，可以花费|Cyan|钻石||购买到|Purple|稀有的道具||，令您事半功倍。'
$ws.Range("E8").Value = 5

$ws.Range("A9").Value = 6
$ws.Range("B9").Value = 20
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = '|您可以使用|Red|卡片商店||了！在卡片商店中，你可以随机购买到各种|Gold|极品卡片||！'
$ws.Range("E9").Value = 8

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = '|您可以进行|Blue|游戏问答||了。通过游戏问答中回答问题，您可以获得一定的|Cyan|阅历||。'
$ws.Range("E10").Value = 9

$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 50
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = '|您可以使用|Yellow|经验瓶||了。经验瓶可以|Red|将战斗中英雄获得经验转化为其他卡片的经验值|。'
$ws.Range("E11").Value = 6

$ws.Range("A12").Value = 100
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = '|您获得了一个|Yellow|经典卡包||。'
$ws.Range("E12").Value = 7

$ws.Range("A13").Value = 101
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = '|新职业-|Gold|Job||，已经可以使用。可以通过城堡面板，转职按钮来完成|Red|转职||。'
$ws.Range("E13").Value = 11

# Match the row styles used for the rest of the data rows (copy from a
# neighbouring, already-correctly-styled row) so the new rows 12/13 pick
# up the same "data row" formatting as row 11 did before the insert.
$ws.Range("A11:E11").Copy() | Out-Null
$ws.Range("A12:E13").PasteSpecial(-4122) | Out-Null

# Re-apply the literal values/styles for row 6 (PasteSpecial above does
# not touch it) and highlight the new unlock row: Level/Type/Icon are
# shown in bold to draw attention to the newly unlocked feature.
$ws.Range("B6").Font.Bold = $true
$ws.Range("C6").Font.Bold = $true
$ws.Range("E6").Font.Bold = $true

# Grow the table so the two new rows become part of "表1" (A3:E13),
# matching Excel's behaviour of auto-expanding a table when new rows are
# appended directly below/within it.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:E13"))

# Leave the selection where Excel would after typing the last new cell.
$ws.Range("D13").Select()
